$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (style) of the last existing data row (251) down to the
# four new rows (252:255) so column A keeps its date style (s="2") and the
# other columns keep the plain numeric style used throughout the table.
$ws.Range("A251:D251").Copy()
$ws.Range("A252:D255").PasteSpecial(-4122)

$data = @(
    @(252, 44326, 0, 17, 99.21213889699445),
    @(253, 44327, 1, 17, 99.21213889699445),
    @(254, 44328, 0, 14, 81.70411438576014),
    @(255, 44329, 1, 8, 46.68806536329151)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
}
